$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (FAPs -> Ccl17 -> Ackr2 -> FAPs) with refreshed TPM values ---
$ws.Range("G2").Value = 0.2999813333333333
$ws.Range("H2").Value = 0.899944
$ws.Range("I2").Value = 0.1904846973062729
$ws.Range("J2").Value = 0.1904846973062729
$ws.Range("Q2").Value = 0.08238307362311111
$ws.Range("R2").Value = 0.741447662608
$ws.Range("S2").Value = 0.1904846973062729
$ws.Range("T2").Value = 0.1904846973062729

# --- Update existing row 3 (MuSCs -> Ccl17 -> Ackr2 -> FAPs) with refreshed TPM values ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2731446666666666
$ws.Range("H3").Value = 0.819434
$ws.Range("I3").Value = 0.1734437225565907
$ws.Range("J3").Value = 0.1734437225565907
$ws.Range("Q3").Value = 0.07501299142088888
$ws.Range("R3").Value = 0.675116922788
$ws.Range("S3").Value = 0.1734437225565907
$ws.Range("T3").Value = 0.1734437225565907

# --- Append new row 4 for the Resolving-Mac sending cluster ---
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Ccl17"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.001705666666667
$ws.Range("H4").Value = 3.005117
$ws.Range("I4").Value = 0.6360715801371364
$ws.Range("J4").Value = 0.6360715801371364
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2746273333333333
$ws.Range("N4").Value = 0.823882
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.2750957560215556
$ws.Range("R4").Value = 2.475861804194
$ws.Range("S4").Value = 0.6360715801371364
$ws.Range("T4").Value = 0.6360715801371364
